$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.261.93'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.656.69'
$ws.Range("D5").Value = '219.38'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '0.5239'
$ws.Range("E6").Value = '  -1.88%  '
$ws.Range("E7").Value = '  -0.63%  '
$ws.Range("D8").Value = '0.2671'
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '0.06370'
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("D10").Value = '20.70'
$ws.Range("D11").Value = '0.07704'
$ws.Range("D12").Value = '4.609'
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("D13").Value = '1.572.60'
$ws.Range("E13").Value = '  -6.11%  '
$ws.Range("D14").Value = '1.885.84'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '0.5651'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").Value = '0.0₅8286'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '65.48'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '26.259.66'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = '4.697'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = '193.03'
$ws.Range("E22").Value = '  -2.58%  '
$ws.Range("D23").Value = '6.010'
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '143.33'
$ws.Range("D26").Value = '0.1202'
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("D27").Value = '7.295'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").Value = '15.93'
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("D29").Value = '1.501'
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = '0.05648'
$ws.Range("E30").Value = '  -4.73%  '
$ws.Range("D31").Value = '1.275'
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").Value = '3.510'
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("D33").Value = '3.356'
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").Value = '1.584'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("D35").Value = '2.807'
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("D36").Value = '0.9481'
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("D37").Value = '2.415'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").Value = '0.5776'
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").Value = '0.01602'
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D41").Value = '2.572'
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("D42").Value = '0.8459'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '101.83'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.018.11'
$ws.Range("E45").Value = '  -5.72%  '
$ws.Range("D46").Value = '1.796.49'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").Value = '58.40'
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").Value = '0.0₈107'
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("E50").Value = '  +3.05%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '8.029'
$ws.Range("E51").Value = '  -0.09%  '
